# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2..16) the existing error-series values in columns
# B..K get shifted one column to the right (B->C, C->D, ... J->K, the old
# K value is dropped) and a newly-computed value is inserted into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert into column B for each row.
$newB = @{
    2  = 0.3648791949059138
    3  = -0.2352699264540507
    4  = -0.05148746350304451
    5  = -0.1333319740152609
    6  = 1.614150253737389
    7  = 0.5701030647716323
    8  = 0.2202779152847414
    9  = 0.5040960054549828
    10 = 0.420735823599318
    11 = -0.1252583916527783
    12 = 0.08824118641116785
    13 = -0.1133200159455487
    14 = 0.1743923273248104
    15 = -0.4559694969238889
    16 = 0.1808172637304477
}

for ($row = 2; $row -le 16; $row++) {

    # Snapshot the current B..K values (columns 2..11) for this row before
    # writing anything back, so the shift doesn't clobber its own source.
    $existing = @()
    for ($col = 2; $col -le 11; $col++) {
        $existing += , $ws.Cells.Item($row, $col).Value()
    }

    # Shift every existing value one column to the right: new column c
    # (3..11, i.e. C..K) gets the old value from column c-1 (2..10, i.e.
    # B..J). Walk from the rightmost column down so column K (index 11,
    # the former J) is written before being needed as a source elsewhere.
    for ($col = 11; $col -ge 3; $col--) {
        $ws.Cells.Item($row, $col).Value = $existing[$col - 3]
    }

    # Insert the newly computed value into column B.
    $ws.Cells.Item($row, 2).Value = $newB[$row]
}

Write-Output "done"
